$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List")

# Copy formatting (including the Hyperlink cell style) from the row above
# so the new row matches the existing pattern (s="5" on column B).
$ws.Range("B50").Copy() | Out-Null
$ws.Range("B51").PasteSpecial(-4122) | Out-Null

# New entry: "Isomorphic Strings"
$ws.Cells.Item(51, 2).Value2 = "Isomorphic Strings"
$ws.Cells.Item(51, 3).Value2 = 1
$ws.Cells.Item(51, 4).Value2 = 1
$ws.Cells.Item(51, 5).Value2 = 46
$ws.Cells.Item(51, 6).Value2 = 0.4
$ws.Cells.Item(51, 7).Value2 = 16.6
$ws.Cells.Item(51, 8).Value2 = 0.42
$ws.Cells.Item(51, 9).Value2 = "https://leetcode.com/problems/isomorphic-strings/submissions/1087078218/"

# Hyperlink for the problem name cell
$ws.Hyperlinks.Add($ws.Range("B51"), "https://leetcode.com/problems/isomorphic-strings/", "", "", "Isomorphic Strings") | Out-Null

# Re-apply the hyperlink-cell formatting (Hyperlinks.Add can reset the style)
$ws.Range("B50").Copy() | Out-Null
$ws.Range("B51").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the active selection as in the authored edit
$ws.Range("G55").Select() | Out-Null
